$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.19701361656189
$ws.Range("B1").Value = 2.045663595199585
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.076266765594482
$ws.Range("E1").Value = 1.210805892944336
